# Fixed update to excel issue
# Shift the Week_Start_Date values forward by one week on the
# "Forecast Comparison" sheet, update the handful of forecast
# values that moved along with them, then refresh the dependent
# Summary sheet figures.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

function Set-TextCell {
    param(
        $Sheet,
        [string]$Address,
        [string]$Text
    )
    # Every value in these columns is stored as plain text in the workbook
    # (even the date-shaped and number-shaped ones). Writing it with a
    # leading apostrophe stops Excel's COM layer from auto-converting
    # date-/number-shaped text into a real date serial or numeric value;
    # resetting the style back to Normal afterwards drops the "stored as
    # text" quote-prefix formatting flag that the apostrophe trick leaves
    # behind, so the cell ends up identical to a plain text cell.
    $cell = $Sheet.Range($Address)
    $cell.Value = "'" + $Text
    $cell.Style = "Normal"
}

# --- Forecast Comparison sheet -------------------------------------------------
# Column B = Week_Start_Date, column D = MyForecast, E = Amazon Mean Forecast,
# F = Amazon P70 Forecast, G = Amazon P80 Forecast, H = Amazon P90 Forecast.

Set-TextCell $wsForecast "B2"  "2025-02-02"
Set-TextCell $wsForecast "B3"  "2025-02-09"

Set-TextCell $wsForecast "B4"  "2025-02-16"
$wsForecast.Range("G4").Value  = 2

Set-TextCell $wsForecast "B5"  "2025-02-23"
$wsForecast.Range("E5").Value  = 1
$wsForecast.Range("G5").Value  = 2

Set-TextCell $wsForecast "B6"  "2025-03-02"

Set-TextCell $wsForecast "B7"  "2025-03-09"
$wsForecast.Range("G7").Value  = 2

Set-TextCell $wsForecast "B8"  "2025-03-16"
$wsForecast.Range("H8").Value  = 4

Set-TextCell $wsForecast "B9"  "2025-03-23"
Set-TextCell $wsForecast "B10" "2025-03-30"
Set-TextCell $wsForecast "B11" "2025-04-06"
Set-TextCell $wsForecast "B12" "2025-04-13"
Set-TextCell $wsForecast "B13" "2025-04-20"

Set-TextCell $wsForecast "B14" "2025-04-27"
$wsForecast.Range("H14").Value = 4

Set-TextCell $wsForecast "B15" "2025-05-04"
$wsForecast.Range("G15").Value = 2
$wsForecast.Range("H15").Value = 4

Set-TextCell $wsForecast "B16" "2025-05-11"
$wsForecast.Range("D16").Value = 0

Set-TextCell $wsForecast "B17" "2025-05-18"

# --- Summary sheet --------------------------------------------------------------
# Every value in Summary!B is plain text too (even the numeric-looking ones).
Set-TextCell $wsSummary "B2"  "2022-12-25 to 2025-01-26"
Set-TextCell $wsSummary "B5"  "11"
Set-TextCell $wsSummary "B12" "0"
Set-TextCell $wsSummary "B13" "2025-02-02"
Set-TextCell $wsSummary "B15" "2025-02-16"
